# Update countries & provincias Spain
# - Swap Nigeria/Finlandia ordering: Finlandia moves up to row 61 (with updated
#   figures), Nigeria moves down to row 62 (keeping its previous figures).
# - Refresh the "datos actualizados" timestamp string.
# - Refresh case/death/recovered counters for several countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 14:35"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1571328
$ws.Range("C4").Value = 745
$ws.Range("E4").Value = 1116540
$ws.Range("G4").Value = 28
$ws.Range("H4").Value = 93561

# India (row 14)
$ws.Range("B14").Value = 107789
$ws.Range("C14").Value = 1314
$ws.Range("D14").Value = 42914
$ws.Range("E14").Value = 61559
$ws.Range("G14").Value = 14
$ws.Range("H14").Value = 3316

# Paises Bajos (row 23)
$ws.Range("B23").Value = 44447
$ws.Range("C23").Value = 198
$ws.Range("G23").Value = 33
$ws.Range("H23").Value = 5748

# Suecia (row 27)
$ws.Range("B27").Value = 31523
$ws.Range("C27").Value = 724
$ws.Range("E27").Value = 22721
$ws.Range("G27").Value = 88
$ws.Range("H27").Value = 3831

# Dinamarca (row 47)
$ws.Range("D47").Value = 9536
$ws.Range("E47").Value = 1027
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 554

# Nigeria / Finlandia reorder (rows 61-62):
# Finlandia takes row 61 with refreshed figures, Nigeria drops to row 62
# retaining its previous (unchanged) figures.
$ws.Range("A61").Value = "Finlandia"
$ws.Range("B61").Value = 6443
$ws.Range("C61").Value = 44
$ws.Range("D61").Value = 5000
$ws.Range("E61").Value = 1142
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 301

$ws.Range("A62").Value = "Nigeria"
$ws.Range("B62").Value = 6401
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 1734
$ws.Range("E62").Value = 4475
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 192

# Croacia (row 82)
$ws.Range("B82").Value = 2234
$ws.Range("C82").Value = 2
$ws.Range("D82").Value = 1978
$ws.Range("E82").Value = 160

# Burkina Faso (row 117)
$ws.Range("B117").Value = 809
$ws.Range("C117").Value = 13
$ws.Range("D117").Value = 661
$ws.Range("E117").Value = 96
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 52
